$wb = $excel.ActiveWorkbook

$wsGlobal = $wb.Worksheets.Item("global_settings")
$wsTech   = $wb.Worksheets.Item("tech")

# [RULES] languages other than english disabled.
# Row 5 (en-US) stays enabled; rows 6-15 (fr-FR, it-IT, de-DE, es-ES,
# lang_brazilian, lang_russian, zh-CN, ja-JP, ko-KR, zh-TW) get their
# [android] (F) and [iOS] (G) flags cleared.
$wsTech.Range("F6:G15").Value = $false

# The edit leaves the "tech" sheet as the active / selected tab, with the
# "global_settings" sheet scrolled/selected at F20 and "tech" selected at G15.
[void]$wsGlobal.Activate()
[void]$wsGlobal.Range("F20").Select()

[void]$wsTech.Activate()
[void]$wsTech.Range("G15").Select()
